# Update "想去人数" (interest count) / "最低票价" (min price) figures to the
# latest scrape snapshot across all four sheets (展览, 演出, 本地生活, 全部类型).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 379
$ws1.Range("F4").Value = 427
$ws1.Range("F8").Value = 1070
$ws1.Range("F9").Value = 1653
$ws1.Range("F10").Value = 6148
$ws1.Range("G10").Value = 68
$ws1.Range("F11").Value = 119
$ws1.Range("F12").Value = 1783
$ws1.Range("F16").Value = 6386
$ws1.Range("F17").Value = 6386
$ws1.Range("F20").Value = 157
$ws1.Range("F21").Value = 101
$ws1.Range("F22").Value = 1680
$ws1.Range("F27").Value = 1476
$ws1.Range("F29").Value = 288
$ws1.Range("F32").Value = 40
$ws1.Range("F34").Value = 3876

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 322
$ws2.Range("F8").Value = 399

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9518
$ws3.Range("F5").Value = 216

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9518
$ws4.Range("F5").Value = 379
$ws4.Range("F6").Value = 427
$ws4.Range("F11").Value = 322
$ws4.Range("F12").Value = 216
$ws4.Range("F13").Value = 1653
$ws4.Range("F14").Value = 6148
$ws4.Range("G14").Value = 68
$ws4.Range("F15").Value = 119
$ws4.Range("F16").Value = 1783
$ws4.Range("F22").Value = 6387
$ws4.Range("F23").Value = 6387
$ws4.Range("F26").Value = 157
$ws4.Range("F27").Value = 101
$ws4.Range("F28").Value = 1680
$ws4.Range("F33").Value = 1476
$ws4.Range("F36").Value = 288
$ws4.Range("F45").Value = 3876
